$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.221.17'
$ws.Range("E2").Value = '  +1.77%  '

# Row 3
$ws.Range("D3").Value = '1.905.54'
$ws.Range("E3").Value = '  +1.79%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.52'
$ws.Range("E5").Value = '  +0.95%  '

# Row 6
$ws.Range("E6").Value = '  +0.06%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4636'
$ws.Range("E7").Value = '  +0.24%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3953'
$ws.Range("E8").Value = '  +2.38%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.71'
$ws.Range("E9").Value = '  +1.05%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07963'
$ws.Range("E10").Value = '  +1.11%  '

# Row 11
$ws.Range("E11").Value = '  +2.48%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.27'
$ws.Range("E12").Value = '  +2.00%  '

# Row 13
$ws.Range("D13").Value = '1.943.87'
$ws.Range("E13").Value = '  +6.01%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.121'
$ws.Range("E14").Value = '  +1.41%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.762'
$ws.Range("E15").Value = '  +1.02%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06957'
$ws.Range("E16").Value = '  +0.06%  '

# Row 17
$ws.Range("E17").Value = '  +0.30%  '

# Row 18
$ws.Range("E18").Value = '  +0.05%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001008'
$ws.Range("E19").Value = '  +0.44%  '

# Row 20
$ws.Range("E20").Value = '  +2.10%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.003'
$ws.Range("E21").Value = '  +0.04%  '

# Row 22
$ws.Range("D22").Value = '29.250.41'
$ws.Range("E22").Value = '  +1.92%  '

# Row 23
$ws.Range("E23").Value = '  +1.55%  '

# Row 24
$ws.Range("E24").Value = '  -0.02%  '

# Row 25
$ws.Range("D25").Value = '2.162.39'
$ws.Range("E25").Value = '  +4.48%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.050'
$ws.Range("E26").Value = '  -2.44%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.95'
$ws.Range("E27").Value = '  +2.59%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.49'
$ws.Range("E28").Value = '  +0.99%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.917'
$ws.Range("E29").Value = '  +0.54%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.995'
$ws.Range("E30").Value = '  +0.42%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '119.13'
$ws.Range("E31").Value = '  -0.13%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09394'
$ws.Range("E32").Value = '  +0.67%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9260'
$ws.Range("E33").Value = '  +0.79%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.346'
$ws.Range("E34").Value = '  +1.06%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.349'
$ws.Range("E35").Value = '  +1.06%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.264'
$ws.Range("E36").Value = '  -1.82%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05834'

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.178'
$ws.Range("E38").Value = '  +2.03%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02105'
$ws.Range("E39").Value = '  +1.34%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.996'
$ws.Range("E40").Value = '  +4.31%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5756'
$ws.Range("E41").Value = '  +2.33%  '

# Row 42
$ws.Range("E42").Value = '  +1.14%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.968'
$ws.Range("E43").Value = '  +1.96%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.04'
$ws.Range("E44").Value = '  +2.31%  '

# Row 45
$ws.Range("E45").Value = '  +2.66%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.213'
$ws.Range("E46").Value = '  +3.14%  '

# Row 47
$ws.Range("E47").Value = '  -1.63%  '

# Row 48
$ws.Range("E48").Value = '  +2.14%  '

# Row 49
$ws.Range("E49").Value = '  +7.30%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '112.18'
$ws.Range("E50").Value = '  -0.60%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.064'
$ws.Range("E51").Value = '  -4.76%  '
